$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header that used to read "cuenta_banco" is renamed to "banco"
# (column P, the 16th column of the header row).
$ws.Range("P1").Value = "banco"

# The whole header row's style is switched from the default "General"
# number format to a plain Text format (numFmtId 49).
$ws.Range("A1:V1").NumberFormat = "@"
